# Add a new row (#24) to the "practice" worksheet for the
# "balanced binary tree" LeetCode question, matching the formatting of the
# most recently added rows (e.g. row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last similarly-shaped row (21) down onto the new
# row 24 so the new cells pick up the same styles (center alignment, wrapped
# header cell, date number format, etc.) as the rest of the table.
$srcFmt = $ws.Range("A21:J21")
$dstFmt = $ws.Range("A24:J24")
$srcFmt.Copy()
$dstFmt.PasteSpecial(-4122)   # xlPasteFormats

# Row 21 (and similar rows) render a bit taller than the sheet default, so
# match that same explicit row height for the new row.
$ws.Rows.Item(24).RowHeight = 17

# Fill in the new question's data.
$ws.Range("A24").Value = 110
$ws.Range("B24").Value = "balanced binary tree"
$ws.Range("C24").Value = "easy"
$ws.Range("D24").Value = "10/1/2023"
$ws.Range("E24").Value = "solved"
$ws.Range("F24").Value = 7
$ws.Range("G24").Value = "tree"
$ws.Range("H24").Value = "dfs/recursive"
$ws.Range("I24").Value = "O(n)"
$ws.Range("J24").Value = "O(1)"

# Update the sheet's remembered selection, as Excel does when a user leaves
# the cursor somewhere after editing.
[void]$ws.Range("I36").Select()
